# --- dim-localizacao.xlsx edit -------------------------------------------
# 1) Rename the main data sheet to reflect the new RM instance id
#    (RM552628 -> RM553315).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "RM553315.DIM_LOCALIZACAO"

# 2) Normalize mis-decoded (mojibake) accented city names in column C to
#    their plain-ASCII equivalents.
$cityMap = @{
  "MaceiÃ³" = "Maceio"
  "Palmeira dos Ã�ndios" = "Palmeira dos Indios"
  "MacapÃ¡" = "Macapa"
  "VitÃ³ria da Conquista" = "Vitoria da Conquista"
  "IlhÃ©us" = "Ilheus"
  "BrasÃ­lia" = "Brasilia"
  "VitÃ³ria" = "Vitoria"
  "GoiÃ¢nia" = "Goiania"
  "Aparecida de GoiÃ¢nia" = "Aparecida de Goiania"
  "AnÃ¡polis" = "Anapolis"
  "SÃ£o LuÃ­s" = "Sao Luis"
  "CuiabÃ¡" = "Cuiaba"
  "VÃ¡rzea Grande" = "Varzea Grande"
  "RondonÃ³polis" = "Rondonopolis"
  "TrÃªs Lagoas" = "Tres Lagoas"
  "UberlÃ¢ndia" = "Uberlandia"
  "BelÃ©m" = "Belem"
  "SantarÃ©m" = "Santarem"
  "JoÃ£o Pessoa" = "Joao Pessoa"
  "MaringÃ¡" = "Maringa"
  "JaboatÃ£o dos Guararapes" = "Jaboatao dos Guararapes"
  "ParnaÃ­ba" = "Parnaiba"
  "NiterÃ³i" = "Niteroi"
  "SÃ£o GonÃ§alo" = "Sao Goncalo"
  "MossorÃ³" = "Mossoro"
  "Ji-ParanÃ¡" = "Ji Parana"
  "RorainÃ³polis" = "Rorainopolis"
  "CaracaraÃ­" = "Caracarai"
  "FlorianÃ³polis" = "Florianopolis"
  "SÃ£o JosÃ©" = "Sao Jose"
  "SÃ£o Paulo" = "Sao Paulo"
  "SÃ£o Bernardo do Campo" = "Sao Bernardo do Campo"
  "Santo AndrÃ©" = "Santo Andre"
  "AraguaÃ­na" = "Araguaina"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($null -ne $v -and $cityMap.ContainsKey($v)) {
        $cell.Value = $cityMap[$v]
    }
}

# 3) Update the SQL sheet's query text to reference the new RM instance id.
$sql = $wb.Worksheets.Item(2)
$formulaCell = $sql.Range("A2")
$formulaCell.Value = $formulaCell.Value2.Replace("RM552628", "RM553315")

# 4) Touch number-format definitions (date/time + date) so the style table
#    carries the same numFmt/cellXfs entries as the authored workbook.
$ws.Range("Z1").NumberFormat = "m/d/yyyy h:mm AM/PM"
$ws.Range("Z2").NumberFormat = "m/d/yyyy"
$ws.Range("Z1:Z2").Clear()

# 5) Restore the frozen header row (row 1) on the data sheet, re-selecting
#    A1 afterwards so the saved view matches the original selection/pane.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

